$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 - this pushes the existing rows 28-33
# (the weekly Coco price records) down to rows 29-34.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new weekly price record.
$ws.Cells.Item(28, 1).Value = 9
$ws.Cells.Item(28, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(28, 3).Value = "Metropolitana"
$ws.Cells.Item(28, 4).Value = 44452
$ws.Cells.Item(28, 5).Value = 13
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100108
$ws.Cells.Item(28, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(28, 9).Value = 100108007
$ws.Cells.Item(28, 10).Value = "Coco"
$ws.Cells.Item(28, 11).Value = "Sin especificar"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 35
$ws.Cells.Item(28, 14).Value = 21000
$ws.Cells.Item(28, 15).Value = 22000
$ws.Cells.Item(28, 16).Value = 21429
$ws.Cells.Item(28, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(28, 18).Value = "Perú"
$ws.Cells.Item(28, 19).Value = 1071
$ws.Cells.Item(28, 20).Value = 20
